$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "F3"  = 2
    "F14" = 3
    "F15" = 2
    "F16" = 3
    "F19" = -1
    "F21" = 1
    "F22" = 0
    "F24" = 0
    "F33" = 0
    "F41" = -1
    "F46" = -1
    "F57" = 0
    "F59" = -1
    "F63" = -1
    "F72" = -2
    "F75" = 11
    "F76" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
